$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format/border/font/alignment) from A52 to new A53
# so the inserted row keeps the same per-column look as the rest of the table.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)

# Update the forecast data table (columns A:E) for rows 2-53.
# A new earliest observation was inserted at the top (old row 2 data shifted down
# by one row), and the y_0_forecast / y_1_forecast values were recomputed for every row.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 0.4235526809466261
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 0.6439341879002525
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = -0.3623658873974311
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 0.1825419310453658
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = -0.5718076928962645
$ws.Range("D4").Value = 2009
$ws.Range("E4").Value = -0.1800933741311961
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("C5").Value = -0.009261555895478946
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = 0.1145211022186787
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = 0.3486139762225005
$ws.Range("D6").Value = 2010
$ws.Range("E6").Value = 0.1555182634501051
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = -1.404263945418582
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = -0.807808220045203
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = -0.1384957661262898
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 0.6938817570587785
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 1.692932643509848
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 0.6262577107155831
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 1.566479473280147
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 0.9614071719361794
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 1.020829760720643
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 1.148272834981245
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 0.7307568962936939
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 1.09290550768979
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = 0.6772121200332215
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 1.258913537332895
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.818818812164257
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 0.9607602172681418
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 1.019715257608911
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 0.9536145745415947
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 0.9180054319587239
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 1.375398114243209
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 2.173959184500385
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 1.566646323486043
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 1.984684278296656
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 1.473274087935805
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 1.707434489469994
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 1.30258347990615
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 1.755995812646982
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 1.681032827388362
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 1.456988786619839
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 1.842797144428188
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 1.946965557828384
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 1.755491062323111
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 1.131202984360957
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 1.657737120813452
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 1.241332692055597
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 1.58004210678635
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 1.260396653238383
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 1.567743002885069
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 1.06432145354225
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 0.776718238020746
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 0.5757500748109434
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 1.030688008679626
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 1.592885137608979
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 1.604795846351514
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 1.308235387832934
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 1.242807488305719
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 1.361817904277696
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 1.316199564471554
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 0.9437384066259158
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 0.904959070968947
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -2.015335584265165
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -1.215549235925828
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -5.210209911466245
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -2.349089443609143
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -4.352425014431304
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 0.03547044462246518
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -4.454337270215236
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = -3.012953608516933
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = -3.579597300369253
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = -1.403103901755631
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = -1.897775264882628
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 4.997412520017441
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = -1.761645650979182
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 3.765721202592909
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 3.304925622412869
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = -0.4653479251390702
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 4.461954539041502
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 0.7797949948739058
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 4.787836378515364
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 1.112953228460167
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 5.20787683103745
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = 3.217995704408838
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.5032298616872488
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 0.34496349151385
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -1.305206755692701
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.5821000732047832
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -0.49899188013105
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 2.610266500707703
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -0.9008525709169546
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.6027009207580036
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = 0.2229020320597241
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = -0.08457672677967265
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = 0.0845726262934221
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.1341520870597357
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = 0.5084754301873051
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.0148008406940292
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = 0.2738544794132824
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 0.2681899963140832
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 0.2312068876759277
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = -0.4074819591325718
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.9724700385226326
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 0.6236501628417823
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 0.9584581489103794
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 0.6809779381435677
